$wb = $excel.ActiveWorkbook

# Hunk 0: ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2140.04
$ws.Cells.Item(40, 10).Value = 2583.6667
$ws.Cells.Item(40, 12).Value = 2583.6667
$ws.Cells.Item(40, 14).Value = -2933.6667

# Hunk 1: ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1779.875
$ws.Cells.Item(98, 9).Value = 1779.875
$ws.Cells.Item(98, 11).Value = 1779.875
$ws.Cells.Item(98, 13).Value = -281.875

# Hunk 2: ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1779.875
$ws.Cells.Item(122, 9).Value = 1779.875
$ws.Cells.Item(122, 11).Value = 5339.625
$ws.Cells.Item(122, 13).Value = -2889.625

# Hunk 3: ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1407.8422
$ws.Cells.Item(135, 10).Value = 2548.3333
$ws.Cells.Item(135, 12).Value = 22934.9997
$ws.Cells.Item(135, 14).Value = -28004.9997

# Hunk 4: ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2080.873
$ws.Cells.Item(138, 9).Value = 1890.05
$ws.Cells.Item(138, 10).Value = 2169.628
$ws.Cells.Item(138, 11).Value = 5670.15
$ws.Cells.Item(138, 12).Value = 6508.884
$ws.Cells.Item(138, 13).Value = -530.1499999999996
$ws.Cells.Item(138, 14).Value = -16788.884

# Hunk 5: ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 5838.2144
$ws.Cells.Item(141, 9).Value = 2185.4
$ws.Cells.Item(141, 11).Value = 6556.200000000001
$ws.Cells.Item(141, 13).Value = -1376.200000000001

# Hunk 6: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1010.7692
$ws.Cells.Item(61, 9).Value = 831.1818
$ws.Cells.Item(61, 10).Value = 1998.5
$ws.Cells.Item(61, 11).Value = 831.1818
$ws.Cells.Item(61, 12).Value = 1998.5
$ws.Cells.Item(61, 13).Value = -619.1818
$ws.Cells.Item(61, 14).Value = -2422.5

# Hunk 7: ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1025
$ws.Cells.Item(63, 9).Value = 2600
$ws.Cells.Item(63, 11).Value = 2600
$ws.Cells.Item(63, 13).Value = -1914

# Hunk 8: ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 1025
$ws.Cells.Item(66, 9).Value = 2600
$ws.Cells.Item(66, 11).Value = 13000
$ws.Cells.Item(66, 13).Value = -9568

# Hunk 9: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 414.54544
$ws.Cells.Item(74, 9).Value = 365.03125
$ws.Cells.Item(74, 11).Value = 365.03125
$ws.Cells.Item(74, 13).Value = 508.96875

# Hunk 10: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 414.54544
$ws.Cells.Item(77, 9).Value = 365.03125
$ws.Cells.Item(77, 11).Value = 1825.15625
$ws.Cells.Item(77, 13).Value = 2542.84375

# Hunk 11: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1010.7692
$ws.Cells.Item(136, 9).Value = 831.1818
$ws.Cells.Item(136, 10).Value = 1998.5
$ws.Cells.Item(136, 11).Value = 2493.5454
$ws.Cells.Item(136, 12).Value = 5995.5
$ws.Cells.Item(136, 13).Value = 56.45460000000003
$ws.Cells.Item(136, 14).Value = -11095.5

# Hunk 12: BSM row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 20000
$ws.Cells.Item(92, 9).Value = 20000
$ws.Cells.Item(92, 11).Value = 20000
$ws.Cells.Item(92, 13).Value = -17504

# Hunk 13: BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 975.4666999999999
$ws.Cells.Item(99, 9).Value = 824.6667
$ws.Cells.Item(99, 11).Value = 824.6667
$ws.Cells.Item(99, 13).Value = 673.3333

# Hunk 14: BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2123.7693
$ws.Cells.Item(107, 9).Value = 1400
$ws.Cells.Item(107, 11).Value = 1400
$ws.Cells.Item(107, 13).Value = 520

# Hunk 15: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1518.5
$ws.Cells.Item(134, 9).Value = 1115.6666
$ws.Cells.Item(134, 10).Value = 3532.6667
$ws.Cells.Item(134, 11).Value = 3346.9998
$ws.Cells.Item(134, 12).Value = 10598.0001
$ws.Cells.Item(134, 13).Value = -811.9998000000001
$ws.Cells.Item(134, 14).Value = -15668.0001

# Hunk 16: BSM row 138
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

# Hunk 17: CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 225.27272
$ws.Cells.Item(12, 10).Value = 218.6
$ws.Cells.Item(12, 12).Value = 655.8
$ws.Cells.Item(12, 14).Value = -1001.8

# Hunk 18: CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 6533.3335
$ws.Cells.Item(55, 9).Value = 100
$ws.Cells.Item(55, 11).Value = 300
$ws.Cells.Item(55, 13).Value = -123

# Hunk 19: CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 6398.8
$ws.Cells.Item(132, 9).Value = 4994
$ws.Cells.Item(132, 10).Value = 6750
$ws.Cells.Item(132, 11).Value = 44946
$ws.Cells.Item(132, 12).Value = 60750
$ws.Cells.Item(132, 13).Value = -42416
$ws.Cells.Item(132, 14).Value = -65810

# Hunk 20: GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 25039688
$ws.Cells.Item(113, 9).Value = 50029800
$ws.Cells.Item(113, 11).Value = 50029800
$ws.Cells.Item(113, 13).Value = -50027630

# Hunk 21: GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 68937.60000000001
$ws.Cells.Item(122, 9).Value = 1988
$ws.Cells.Item(122, 11).Value = 5964
$ws.Cells.Item(122, 13).Value = -3514

# Hunk 22: GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value = 63995.668
$ws.Cells.Item(134, 10).Value = 63995.668
$ws.Cells.Item(134, 12).Value = 191987.004
$ws.Cells.Item(134, 14).Value = -197057.004

# Hunk 23: GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 210000
$ws.Cells.Item(136, 10).Value = 210000
$ws.Cells.Item(136, 12).Value = 630000
$ws.Cells.Item(136, 14).Value = -635100

# Hunk 24: LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 14).ClearContents()

# Hunk 25: LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1057.9166
$ws.Cells.Item(16, 9).Value = 1057.9166
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1057.9166
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -887.9166
$ws.Cells.Item(16, 14).ClearContents()

# Hunk 26: LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2366.6667
$ws.Cells.Item(40, 9).Value = 2366.6667
$ws.Cells.Item(40, 11).Value = 2366.6667
$ws.Cells.Item(40, 13).Value = -2230.6667

# Hunk 27: LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 955.35297
$ws.Cells.Item(93, 9).Value = 864.9375
$ws.Cells.Item(93, 11).Value = 864.9375
$ws.Cells.Item(93, 13).Value = 383.0625

# Hunk 28: LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3165.5
$ws.Cells.Item(122, 9).Value = 2897.6667
$ws.Cells.Item(122, 11).Value = 8693.000100000001
$ws.Cells.Item(122, 13).Value = -6243.000100000001

# Hunk 29: LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 14).ClearContents()

# Hunk 30: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2666
$ws.Cells.Item(132, 9).Value = 2499.5
$ws.Cells.Item(132, 10).Value = 2999
$ws.Cells.Item(132, 11).Value = 7498.5
$ws.Cells.Item(132, 12).Value = 8997
$ws.Cells.Item(132, 13).Value = -4968.5
$ws.Cells.Item(132, 14).Value = -14057

# Hunk 31: WVR row 69
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 6590.3335
$ws.Cells.Item(69, 10).Value = 6590.3335
$ws.Cells.Item(69, 12).Value = 6590.3335
$ws.Cells.Item(69, 14).Value = -8088.3335

# Hunk 32: WVR row 72
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(72, 8).Value = 6590.3335
$ws.Cells.Item(72, 10).Value = 6590.3335
$ws.Cells.Item(72, 12).Value = 19771.0005
$ws.Cells.Item(72, 14).Value = -27259.0005

# Hunk 33: WVR row 120
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(120, 8).Value = 15000
$ws.Cells.Item(120, 10).Value = 15000
$ws.Cells.Item(120, 12).Value = 15000
$ws.Cells.Item(120, 14).Value = -24676

# Hunk 34: WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 4331.3335
$ws.Cells.Item(126, 9).Value = 3997.5
$ws.Cells.Item(126, 11).Value = 11992.5
$ws.Cells.Item(126, 13).Value = -9522.5

# Hunk 35: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1691.7333
$ws.Cells.Item(132, 9).Value = 1752
$ws.Cells.Item(132, 10).Value = 1300
$ws.Cells.Item(132, 11).Value = 5256
$ws.Cells.Item(132, 12).Value = 3900
$ws.Cells.Item(132, 13).Value = -2726
$ws.Cells.Item(132, 14).Value = -8960
